# Insert a new weekly record at row 38 ("Fruta / hortaliza, semanal"):
# this pushes the existing rows 38-130 down to 39-131 (dimension grows to
# A1:R131) and fills the freed row 38 with the new observation, reusing
# the constant metadata columns (A,B,C,E,F,G,H,I,N,O,Q,R) from the record
# that is now immediately below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 38:130 down by one row, leaving row 38 blank (keeps formatting).
$ws.Rows.Item(38).Insert()

# Carry over the unchanged metadata from the row that got pushed to 39.
$ws.Range("A38").Value = $ws.Range("A39").Value()
$ws.Range("B38").Value = $ws.Range("B39").Value()
$ws.Range("C38").Value = $ws.Range("C39").Value()
$ws.Range("E38").Value = $ws.Range("E39").Value()
$ws.Range("F38").Value = $ws.Range("F39").Value()
$ws.Range("G38").Value = $ws.Range("G39").Value()
$ws.Range("H38").Value = $ws.Range("H39").Value()
$ws.Range("I38").Value = $ws.Range("I39").Value()
$ws.Range("N38").Value = $ws.Range("N39").Value()
$ws.Range("O38").Value = $ws.Range("O39").Value()
$ws.Range("Q38").Value = $ws.Range("Q39").Value()
$ws.Range("R38").Value = $ws.Range("R39").Value()

# New weekly price observation.
$ws.Range("D38").Value = 45071
$ws.Range("J38").Value = 400
$ws.Range("K38").Value = 22000
$ws.Range("L38").Value = 24000
$ws.Range("M38").Value = 23000
$ws.Range("P38").Value = 920
